$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 4 new rows (8-11) that duplicate existing entries (Moss, Arlo, Wern Joe,
# Kinloch Heath respectively), mirroring the source rows 3, 2, 7 and 6 so that
# values, shared-string reuse, number formatting and date styling all match.
$ws.Range("A3:H3").Copy($ws.Range("A8:H8"))
$ws.Range("A2:H2").Copy($ws.Range("A9:H9"))
$ws.Range("A7:H7").Copy($ws.Range("A10:H10"))
$ws.Range("A6:H6").Copy($ws.Range("A11:H11"))

# Move the active selection to reflect where the user ended up after the edit.
$ws.Range("B21").Select()
